# Insert a new price record as row 62 in the "Puerro" (Vega Modelo de
# Temuco) weekly log, pushing the existing rows 62-138 down to 63-139.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(62).Insert()

$ws.Range("A62").Value = 10
$ws.Range("B62").Value = "Vega Modelo de Temuco"
$ws.Range("C62").Value = "La Araucanía"
$ws.Range("D62").Value = 44483
$ws.Range("E62").Value = 9
$ws.Range("F62").Value = 100112005
$ws.Range("G62").Value = "Puerro"
$ws.Range("H62").Value = "Azul de Maquehue"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 80
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = 6500
$ws.Range("N62").Value = "$/docena de paquetes"
$ws.Range("O62").Value = "Provincia de Cautín"
$ws.Range("P62").Value = 542
$ws.Range("Q62").Value = 12
$ws.Range("R62").Value = "Hortaliza"
